$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Policy_Id" values (duplicated hash "5fd0709530a434204c3007d5") in column D
# for the data rows (2-6) are no longer needed; clear their contents while leaving
# the rest of the row (and the header) untouched.
$ws.Range("D2:D6").ClearContents()
